$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itgb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.021590333333334
$ws.Range("H2").Value = 6.064771
$ws.Range("I2").Value = 0.01116262347650641
$ws.Range("J2").Value = 0.01116262347650641
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.128124
$ws.Range("N2").Value = 0.384372
$ws.Range("O2").Value = 0.3522399658364659
$ws.Range("P2").Value = 0.352239965836466
$ws.Range("Q2").Value = 0.259014239868
$ws.Range("R2").Value = 2.331128158812001
$ws.Range("S2").Value = 0.003931922112009952
$ws.Range("T2").Value = 0.003931922112009952

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itgb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.021590333333334
$ws.Range("H3").Value = 6.064771
$ws.Range("I3").Value = 0.01116262347650641
$ws.Range("J3").Value = 0.01116262347650641
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.2109236666666666
$ws.Range("N3").Value = 0.632771
$ws.Range("O3").Value = 0.5798737562109268
$ws.Range("P3").Value = 0.5798737562109268
$ws.Range("Q3").Value = 0.4264012456045556
$ws.Range("R3").Value = 3.837611210441
$ws.Range("S3").Value = 0.006472912404490048
$ws.Range("T3").Value = 0.006472912404490048

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itgb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.021590333333334
$ws.Range("H4").Value = 6.064771
$ws.Range("I4").Value = 0.01116262347650641
$ws.Range("J4").Value = 0.01116262347650641
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.024693
$ws.Range("N4").Value = 0.07407900000000001
$ws.Range("O4").Value = 0.06788627795260727
$ws.Range("P4").Value = 0.06788627795260727
$ws.Range("Q4").Value = 0.04991913010100001
$ws.Range("R4").Value = 0.4492721709090001
$ws.Range("S4").Value = 0.0007577889600064136
$ws.Range("T4").Value = 0.0007577889600064136

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itgb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 25.140634
$ws.Range("H5").Value = 75.421902
$ws.Range("I5").Value = 0.1388191398995883
$ws.Range("J5").Value = 0.1388191398995883
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.128124
$ws.Range("N5").Value = 0.384372
$ws.Range("O5").Value = 0.3522399658364659
$ws.Range("P5").Value = 0.352239965836466
$ws.Range("Q5").Value = 3.221118590616001
$ws.Range("R5").Value = 28.990067315544
$ws.Range("S5").Value = 0.04889764909567857
$ws.Range("T5").Value = 0.04889764909567856

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itgb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 25.140634
$ws.Range("H6").Value = 75.421902
$ws.Range("I6").Value = 0.1388191398995883
$ws.Range("J6").Value = 0.1388191398995883
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2109236666666666
$ws.Range("N6").Value = 0.632771
$ws.Range("O6").Value = 0.5798737562109268
$ws.Range("P6").Value = 0.5798737562109268
$ws.Range("Q6").Value = 5.302754705604666
$ws.Range("R6").Value = 47.724792350442
$ws.Range("S6").Value = 0.08049757608754442
$ws.Range("T6").Value = 0.0804975760875444

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itgb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 25.140634
$ws.Range("H7").Value = 75.421902
$ws.Range("I7").Value = 0.1388191398995883
$ws.Range("J7").Value = 0.1388191398995883
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.024693
$ws.Range("N7").Value = 0.07407900000000001
$ws.Range("O7").Value = 0.06788627795260727
$ws.Range("P7").Value = 0.06788627795260727
$ws.Range("Q7").Value = 0.6207976753620001
$ws.Range("R7").Value = 5.587179078258001
$ws.Range("S7").Value = 0.009423914716365325
$ws.Range("T7").Value = 0.009423914716365323

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itgb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 153.9412893333333
$ws.Range("H8").Value = 461.823868
$ws.Range("I8").Value = 0.8500182366239053
$ws.Range("J8").Value = 0.8500182366239052
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.128124
$ws.Range("N8").Value = 0.384372
$ws.Range("O8").Value = 0.3522399658364659
$ws.Range("P8").Value = 0.352239965836466
$ws.Range("Q8").Value = 19.723573754544
$ws.Range("R8").Value = 177.512163790896
$ws.Range("S8").Value = 0.2994103946287774
$ws.Range("T8").Value = 0.2994103946287774

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itgb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 153.9412893333333
$ws.Range("H9").Value = 461.823868
$ws.Range("I9").Value = 0.8500182366239053
$ws.Range("J9").Value = 0.8500182366239052
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2109236666666666
$ws.Range("N9").Value = 0.632771
$ws.Range("O9").Value = 0.5798737562109268
$ws.Range("P9").Value = 0.5798737562109268
$ws.Range("Q9").Value = 32.46986119758089
$ws.Range("R9").Value = 292.228750778228
$ws.Range("S9").Value = 0.4929032677188924
$ws.Range("T9").Value = 0.4929032677188923

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itgb6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 153.9412893333333
$ws.Range("H10").Value = 461.823868
$ws.Range("I10").Value = 0.8500182366239053
$ws.Range("J10").Value = 0.8500182366239052
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.024693
$ws.Range("N10").Value = 0.07407900000000001
$ws.Range("O10").Value = 0.06788627795260727
$ws.Range("P10").Value = 0.06788627795260727
$ws.Range("Q10").Value = 3.801272257508001
$ws.Range("R10").Value = 34.211450317572
$ws.Range("S10").Value = 0.05770457427623553
$ws.Range("T10").Value = 0.05770457427623552
